# Append new trading-day rows (2701-2761) to the CFG_sentiment sheet,
# extending the data from row 2700 (date serial 45825) through row 2761
# (date serial 45915), mirroring the existing Date/Sentiment columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(45826,45828,45831,45832,45833,45834,45835,45838,45839,45840,45841,45845,45846,45847,45848,45849,45852,45853,45854,45855,45856,45859,45860,45861,45862,45863,45866,45867,45868,45869,45870,45873,45874,45875,45876,45877,45880,45881,45882,45883,45884,45887,45888,45889,45890,45891,45894,45895,45896,45897,45898,45902,45903,45904,45905,45908,45909,45910,45911,45912,45915)

$row = 2701
foreach ($d in $dates) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $d
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = 0

    $row = $row + 1
}
